$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "56.817.26"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.967.72"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("D5").Value = "'496.39"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").Value = "'136.36"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.425"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "'7.27"
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").Value = "'0.354"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Value = "3.486.64"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").Value = "'25.61"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "'0.0000156"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "56.875.57"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'6.03"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "2.976.54"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "'12.51"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'7.73"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "'317.95"
$ws.Range("E21").Value = "  -2.58%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").Value = "'0.483"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "'63.15"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -6.07%  "
$ws.Range("D28").Value = "0.0₃0881"
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("D29").Value = "'6.51"
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").Value = "'7.07"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").Value = "'1.75"
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  -6.94%  "
$ws.Range("D33").Value = "'20.04"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").Value = "'153.01"
$ws.Range("D35").Value = "'4.57"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'5.70"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'1.24"
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("D38").Value = "'23.79"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'0.0660"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("D40").Value = "2.997.97"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").Value = "'37.42"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "'3.69"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'0.638"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("D45").Value = "2.190.57"
$ws.Range("E45").Value = "  -4.43%  "
$ws.Range("D46").Value = "'1.37"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("D47").Value = "'0.940"
$ws.Range("E47").Value = "  -6.30%  "
$ws.Range("D48").Value = "'5.89"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "'0.0233"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("D50").Value = "'18.94"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").Value = "'1.77"
$ws.Range("E51").Value = "  -8.70%  "
